# Applies the commit: adds a new "Financial Instrument Informatio" sheet
# with security reference data, and trims the "Dividends" sheet down to
# its first block of rows (removing the duplicated rows 8-13 and the
# now-unused "Code" column E).

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Trim the "Dividends" sheet: drop column E ("Code") and the
#    duplicated rows 8-13.
# ------------------------------------------------------------------
$div = $wb.Worksheets.Item("Dividends")
$div.Range("A8:E13").EntireRow.Delete()
$div.Columns.Item(5).Delete()

# ------------------------------------------------------------------
# 2) Add the new "Financial Instrument Informatio" sheet at the end
#    of the workbook, with the header row + 15 data rows.
# ------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$fii = $wb.Worksheets.Add($null, $last)
$fii.Name = "Financial Instrument Informatio"

$headers = @("Asset Category","Symbol","Description","Conid","Security ID","Underlying","Listing Exch","Multiplier","Type","Code")
for ($c = 0; $c -lt $headers.Count; $c++) {
    $cell = $fii.Cells.Item(1, $c + 1)
    $cell.Value = $headers[$c]
}
$headerRow = $fii.Range("A1:J1")
$headerRow.Font.Bold = $true
$headerRow.Borders.LineStyle = 1
$headerRow.HorizontalAlignment = -4108
$headerRow.VerticalAlignment = -4160

# Columns D (Conid) and H (Multiplier) hold numeric-looking text ("314449552",
# "1", ...) that must stay text (the source export never stores true numbers -
# everything is an inline string), so the whole data block is pre-formatted
# as Text before any values are written - this mirrors typing into a
# Text-formatted column in Excel and stops the usual numeric auto-conversion.
# Column J ("Code") has no data in this export (always blank) so it is left
# untouched.
$dataRange = $fii.Range("A2:I16")
$dataRange.NumberFormat = "@"

$rows = @(
    @("Stocks","AEEM","AMUNDI MSCI EMERG MARK","314449552","LU1681045370","AEEM","SBF","1","ETF"),
    @("Stocks","AJG","ARTHUR J GALLAGHER & CO","4325","US3635761097","AJG","NYSE","1","COMMON"),
    @("Stocks","BXMT","BLACKSTONE MORTGAGE TRU-CL A","127149807","US09257W1009","BXMT","NYSE","1","REIT"),
    @("Stocks","CSPX","ISHARES CORE S&P 500","75776072","IE00B5BMR087","SXR8","IBIS2","1","ETF"),
    @("Stocks","CSX","CSX CORP","6150","US1264081035","CSX","NASDAQ","1","COMMON"),
    @("Stocks","EMD","WESTERN ASSET EMRG MRKT DBT","41073515","US95766A1016","EMD","NYSE","1","CLOSED-END FUND"),
    @("Stocks","GLD","SPDR GOLD SHARES","51529211","US78463V1070","GLD","ARCA","1","ETF"),
    @("Stocks","JPC","NUVEEN PREFERED & INCOME OPP","17635192","US67073B1061","JPC","NYSE","1","CLOSED-END FUND"),
    @("Stocks","LRCX","LAM RESEARCH CORP","732440574","US5128073062","LRCX","NASDAQ","1","COMMON"),
    @("Stocks","MO","ALTRIA GROUP INC","9769","US02209S1033","MO","NYSE","1","COMMON"),
    @("Stocks","QRTEP","QURATE RETAIL INC","442948738","US74915M3088","QRTEP","NASDAQ","1","PUBLIC"),
    @("Stocks","RGLD","ROYAL GOLD INC","4817403","US7802871084","RGLD","NASDAQ","1","COMMON"),
    @("Stocks","TLT","ISHARES 20+ YEAR TREASURY BD","15547841","US4642874329","TLT","NASDAQ","1","ETF"),
    @("Stocks","V","VISA INC-CLASS A SHARES","49462172","US92826C8394","V","NYSE","1","COMMON"),
    @("Stocks","XLB","MATERIALS SELECT SECTOR SPDR","4215200","US81369Y1001","XLB","ARCA","1","ETF")
)

$r = 2
foreach ($row in $rows) {
    for ($c = 0; $c -lt $row.Count; $c++) {
        $fii.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r = $r + 1
}

Write-Output "done"
